$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 1 (the "prompt" header row); remaining rows shift up
$ws.Rows.Item(1).Delete()

# Set column A width to match the diff (auto-size/best-fit to contents)
$ws.Columns.Item(1).AutoFit()

# Update the selection to A11 as shown in the diff
$ws.Range("A11").Select()
